# Recalculated Spencer's method force/thrust results for rows 2-21 (columns I-O, Z-AD).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2").Value = [double]"19348.16829700482"
$ws.Range("L2").Value = [double]"3190.812714091454"
$ws.Range("N2").Value = [double]"19083.2473827914"
$ws.Range("AA2").Value = [double]"0.8380083762482684"
$ws.Range("AC2").Value = [double]"35.56498240708125"
$ws.Range("AD2").Value = [double]"5.456968210637569e-11"
$ws.Range("I3").Value = [double]"19348.16829700482"
$ws.Range("J3").Value = [double]"49139.2201214175"
$ws.Range("K3").Value = [double]"3190.812714091454"
$ws.Range("L3").Value = [double]"8103.81871384848"
$ws.Range("M3").Value = [double]"19083.2473827914"
$ws.Range("N3").Value = [double]"48466.39120456774"
$ws.Range("Z3").Value = [double]"9.711840761553574"
$ws.Range("AA3").Value = [double]"4.99192567737439"
$ws.Range("AB3").Value = [double]"35.56498240708125"
$ws.Range("AC3").Value = [double]"30.84506732290206"
$ws.Range("AD3").Value = [double]"-5.093170329928398e-11"
$ws.Range("I4").Value = [double]"49139.2201214175"
$ws.Range("J4").Value = [double]"81902.05792031095"
$ws.Range("K4").Value = [double]"8103.81871384848"
$ws.Range("L4").Value = [double]"13506.91826279175"
$ws.Range("M4").Value = [double]"48466.39120456774"
$ws.Range("N4").Value = [double]"80780.63041734839"
$ws.Range("Z4").Value = [double]"12.1111654373174"
$ws.Range("AA4").Value = [double]"8.607194139985943"
$ws.Range("AB4").Value = [double]"30.84506732290206"
$ws.Range("AC4").Value = [double]"27.3410960255706"
$ws.Range("AD4").Value = [double]"-1.600710675120354e-10"
$ws.Range("I5").Value = [double]"81902.05792031095"
$ws.Range("J5").Value = [double]"113990.1748725663"
$ws.Range("K5").Value = [double]"13506.91826279175"
$ws.Range("L5").Value = [double]"18798.74589064839"
$ws.Range("M5").Value = [double]"80780.63041734839"
$ws.Range("N5").Value = [double]"112429.3872633711"
$ws.Range("Z5").Value = [double]"14.32270204194236"
$ws.Range("AA5").Value = [double]"11.73099440922736"
$ws.Range("AB5").Value = [double]"27.3410960255706"
$ws.Range("AC5").Value = [double]"24.74938839285561"
$ws.Range("AD5").Value = [double]"1.891748979687691e-10"
$ws.Range("I6").Value = [double]"113990.1748725663"
$ws.Range("J6").Value = [double]"143389.3898543167"
$ws.Range("K6").Value = [double]"18798.74589064839"
$ws.Range("L6").Value = [double]"23647.13192430714"
$ws.Range("M6").Value = [double]"112429.3872633711"
$ws.Range("N6").Value = [double]"141426.0593898722"
$ws.Range("Z6").Value = [double]"16.2588280542065"
$ws.Range("AA6").Value = [double]"14.42946989550886"
$ws.Range("AB6").Value = [double]"24.74938839285561"
$ws.Range("AC6").Value = [double]"22.92003023415797"
$ws.Range("AD6").Value = [double]"1.455191522836685e-11"
$ws.Range("I7").Value = [double]"143389.3898543167"
$ws.Range("J7").Value = [double]"168924.4270559396"
$ws.Range("K7").Value = [double]"23647.13192430714"
$ws.Range("L7").Value = [double]"27858.25517416793"
$ws.Range("M7").Value = [double]"141426.0593898722"
$ws.Range("N7").Value = [double]"166611.4632155555"
$ws.Range("Z7").Value = [double]"17.91001219153198"
$ws.Range("AA7").Value = [double]"16.75203145178586"
$ws.Range("AB7").Value = [double]"22.92003023415797"
$ws.Range("AC7").Value = [double]"21.76204949441185"
$ws.Range("AD7").Value = [double]"2.037268131971359e-10"
$ws.Range("I8").Value = [double]"168924.4270559396"
$ws.Range("J8").Value = [double]"191349.1149094407"
$ws.Range("K8").Value = [double]"27858.25517416793"
$ws.Range("L8").Value = [double]"31556.43362776139"
$ws.Range("M8").Value = [double]"166611.4632155555"
$ws.Range("N8").Value = [double]"188729.1055277988"
$ws.Range("Z8").Value = [double]"19.26211568824663"
$ws.Range("AA8").Value = [double]"18.55782516643347"
$ws.Range("AB8").Value = [double]"21.76204949441185"
$ws.Range("AC8").Value = [double]"21.05775897259869"
$ws.Range("AD8").Value = [double]"-2.037268131971359e-10"
$ws.Range("I9").Value = [double]"191349.1149094407"
$ws.Range("J9").Value = [double]"210984.5370118085"
$ws.Range("K9").Value = [double]"31556.43362776139"
$ws.Range("L9").Value = [double]"34794.61894479147"
$ws.Range("M9").Value = [double]"188729.1055277988"
$ws.Range("N9").Value = [double]"208095.6735508309"
$ws.Range("Z9").Value = [double]"20.17849604217865"
$ws.Range("AA9").Value = [double]"19.87355766569718"
$ws.Range("AB9").Value = [double]"21.05775897259869"
$ws.Range("AC9").Value = [double]"20.75282059611722"
$ws.Range("AD9").Value = [double]"-4.656612873077393e-10"
$ws.Range("I10").Value = [double]"210984.5370118085"
$ws.Range("J10").Value = [double]"226458.2183885399"
$ws.Range("K10").Value = [double]"34794.61894479147"
$ws.Range("L10").Value = [double]"37346.468738155"
$ws.Range("M10").Value = [double]"208095.6735508309"
$ws.Range("N10").Value = [double]"223357.4846484927"
$ws.Range("Z10").Value = [double]"20.66333815909925"
$ws.Range("AA10").Value = [double]"20.8448233195438"
$ws.Range("AB10").Value = [double]"20.75282059611722"
$ws.Range("AC10").Value = [double]"20.93430575656177"
$ws.Range("AD10").Value = [double]"2.037268131971359e-10"
$ws.Range("I11").Value = [double]"226458.2183885399"
$ws.Range("J11").Value = [double]"250494.4630372593"
$ws.Range("K11").Value = [double]"37346.468738155"
$ws.Range("L11").Value = [double]"41310.41787519135"
$ws.Range("M11").Value = [double]"223357.4846484927"
$ws.Range("N11").Value = [double]"247064.6178377266"
$ws.Range("O11").Value = [double]"15890.95654767201"
$ws.Range("Z11").Value = [double]"20.75948831431964"
$ws.Range("AA11").Value = [double]"23.51913132582883"
$ws.Range("AB11").Value = [double]"20.93430575656177"
$ws.Range("AC11").Value = [double]"23.69394876807096"
$ws.Range("AD11").Value = [double]"-4.656612873077393e-10"
$ws.Range("I12").Value = [double]"250494.4630372593"
$ws.Range("J12").Value = [double]"261679.5647364778"
$ws.Range("K12").Value = [double]"41310.41787519135"
$ws.Range("L12").Value = [double]"43155.01443660317"
$ws.Range("M12").Value = [double]"247064.6178377266"
$ws.Range("N12").Value = [double]"258096.5697750535"
$ws.Range("O12").Value = [double]"15890.95654767201"
$ws.Range("Z12").Value = [double]"22.20816468751246"
$ws.Range("AA12").Value = [double]"26.28306114656847"
$ws.Range("AB12").Value = [double]"23.69394876807096"
$ws.Range("AC12").Value = [double]"27.76884522712698"
$ws.Range("AD12").Value = [double]"3.492459654808044e-10"
$ws.Range("I13").Value = [double]"261679.5647364778"
$ws.Range("J13").Value = [double]"250661.4303662877"
$ws.Range("K13").Value = [double]"43155.01443660317"
$ws.Range("L13").Value = [double]"41337.95337457939"
$ws.Range("M13").Value = [double]"258096.5697750535"
$ws.Range("N13").Value = [double]"247229.299000087"
$ws.Range("Z13").Value = [double]"23.86213133810137"
$ws.Range("AA13").Value = [double]"26.74633722700864"
$ws.Range("AB13").Value = [double]"27.76884522712698"
$ws.Range("AC13").Value = [double]"30.65305111603425"
$ws.Range("AD13").Value = [double]"-8.731149137020111e-11"
$ws.Range("I14").Value = [double]"250661.4303662877"
$ws.Range("J14").Value = [double]"229856.1039812186"
$ws.Range("K14").Value = [double]"41337.95337457939"
$ws.Range("L14").Value = [double]"37906.83271595985"
$ws.Range("M14").Value = [double]"247229.299000087"
$ws.Range("N14").Value = [double]"226708.8453741252"
$ws.Range("Z14").Value = [double]"23.39132805773983"
$ws.Range("AA14").Value = [double]"27.38570615235"
$ws.Range("AB14").Value = [double]"30.65305111603425"
$ws.Range("AC14").Value = [double]"34.64742921064442"
$ws.Range("I15").Value = [double]"229856.1039812186"
$ws.Range("J15").Value = [double]"200348.4139407361"
$ws.Range("K15").Value = [double]"37906.83271595985"
$ws.Range("L15").Value = [double]"33040.55746450792"
$ws.Range("M15").Value = [double]"226708.8453741252"
$ws.Range("N15").Value = [double]"197605.1834618796"
$ws.Range("Z15").Value = [double]"22.86699718424556"
$ws.Range("AA15").Value = [double]"28.16298619527942"
$ws.Range("AB15").Value = [double]"34.64742921064442"
$ws.Range("AC15").Value = [double]"39.94341822167827"
$ws.Range("AD15").Value = [double]"-5.529727786779404e-10"
$ws.Range("I16").Value = [double]"200348.4139407361"
$ws.Range("J16").Value = [double]"163473.3948854044"
$ws.Range("K16").Value = [double]"33040.55746450792"
$ws.Range("L16").Value = [double]"26959.29551619566"
$ws.Range("M16").Value = [double]"197605.1834618796"
$ws.Range("N16").Value = [double]"161235.0682098339"
$ws.Range("Z16").Value = [double]"22.32661261347818"
$ws.Range("AA16").Value = [double]"29.36126438301683"
$ws.Range("AB16").Value = [double]"39.94341822167827"
$ws.Range("AC16").Value = [double]"46.97806999121693"
$ws.Range("AD16").Value = [double]"-3.492459654808044e-10"
$ws.Range("I17").Value = [double]"163473.3948854044"
$ws.Range("J17").Value = [double]"121374.339180751"
$ws.Range("K17").Value = [double]"26959.29551619566"
$ws.Range("L17").Value = [double]"20016.50898820958"
$ws.Range("M17").Value = [double]"161235.0682098339"
$ws.Range("N17").Value = [double]"119712.4453826289"
$ws.Range("Z17").Value = [double]"21.96464742304845"
$ws.Range("AA17").Value = [double]"31.6904460486558"
$ws.Range("AB17").Value = [double]"46.97806999121693"
$ws.Range("AC17").Value = [double]"56.70386861682428"
$ws.Range("AD17").Value = [double]"1.891748979687691e-10"
$ws.Range("I18").Value = [double]"121374.339180751"
$ws.Range("J18").Value = [double]"76999.1198049528"
$ws.Range("K18").Value = [double]"20016.50898820958"
$ws.Range("L18").Value = [double]"12698.34780616046"
$ws.Range("M18").Value = [double]"119712.4453826289"
$ws.Range("N18").Value = [double]"75944.82479887128"
$ws.Range("Z18").Value = [double]"22.33020269358826"
$ws.Range("AA18").Value = [double]"37.51260777616771"
$ws.Range("AB18").Value = [double]"56.70386861682428"
$ws.Range("AC18").Value = [double]"71.88627369940373"
$ws.Range("AD18").Value = [double]"5.820766091346741e-11"
$ws.Range("I19").Value = [double]"76999.1198049528"
$ws.Range("J19").Value = [double]"43048.61290702577"
$ws.Range("K19").Value = [double]"12698.34780616046"
$ws.Range("L19").Value = [double]"7099.383222183529"
$ws.Range("M19").Value = [double]"75944.82479887128"
$ws.Range("N19").Value = [double]"42459.17840801356"
$ws.Range("Z19").Value = [double]"26.96587681556761"
$ws.Range("AA19").Value = [double]"50.17647007820467"
$ws.Range("AB19").Value = [double]"71.88627369940373"
$ws.Range("AC19").Value = [double]"95.09686696204079"
$ws.Range("AD19").Value = [double]"-2.291926648467779e-10"
$ws.Range("I20").Value = [double]"43048.61290702577"
$ws.Range("J20").Value = [double]"13627.07479317009"
$ws.Range("K20").Value = [double]"7099.383222183529"
$ws.Range("L20").Value = [double]"2247.315758187477"
$ws.Range("M20").Value = [double]"42459.17840801356"
$ws.Range("N20").Value = [double]"13440.48880441682"
$ws.Range("Z20").Value = [double]"38.39334017957948"
$ws.Range("AA20").Value = [double]"124.1855323624137"
$ws.Range("AB20").Value = [double]"95.09686696204079"
$ws.Range("AC20").Value = [double]"180.889059144875"
$ws.Range("AD20").Value = [double]"-2.746673999354243e-10"
$ws.Range("I21").Value = [double]"13627.07479317009"
$ws.Range("K21").Value = [double]"2247.315758187477"
$ws.Range("M21").Value = [double]"13440.48880441682"
$ws.Range("Z21").Value = [double]"38.39334017957948"
$ws.Range("AA21").Value = [double]"124.1855323624137"
$ws.Range("AB21").Value = [double]"180.889059144875"
